# The records held in worksheet rows 2-8 were re-sorted upstream: each row's
# species-observation data moved to a different row while the sheet's other
# (unchanged) columns stay put. Snapshot the columns that actually move
# (Id, Taxonsorteringsordning, Rodlistade, TaxonId, Artnamn, Vetenskapligt
# namn, Auktor, Ost, Nord) for every affected row first, then write them back
# out in their new order so earlier writes never clobber data a later row
# still needs to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 8
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot every source row before writing anything back.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowValues
}

# new row -> old row it should receive its data from.
$mapping = @{
    2 = 4
    3 = 6
    4 = 8
    5 = 2
    6 = 3
    7 = 5
    8 = 7
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $srcValues = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $srcValues[$col]
    }
}
